# Updates the NATMI LR-pair table with recomputed TPM stats and
# adds the new "Resolving-Mac" sending-cluster rows (14-19).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$table = @(
    @("FAPs", "Ereg", "Egfr", "ECs", 3, 1, 0.8489523333333334, 2.546857, 0.6346554277922418, 0.634655427792242, 3, 1, 1.635346666666667, 4.90604, 0.02683720313876748, 0.02683720313876747, 1.388331368475556, 12.49498231628, 0.01703237663878177, 0.01703237663878177),
    @("FAPs", "Ereg", "Egfr", "FAPs", 3, 1, 0.8489523333333334, 2.546857, 0.6346554277922418, 0.634655427792242, 3, 1, 45.91636366666668, 137.749091, 0.7535202194331003, 0.7535202194331002, 38.98080407299857, 350.8272366569871, 0.4782256972144182, 0.4782256972144183),
    @("FAPs", "Ereg", "Egfr", "Inflammatory-Mac", 3, 1, 0.8489523333333334, 2.546857, 0.6346554277922418, 0.634655427792242, 3, 1, 0.8399643333333332, 2.519893, 0.01378441275019327, 0.01378441275019327, 0.713089680700111, 6.417807126301, 0.00874835237083874, 0.008748352370838743),
    @("FAPs", "Ereg", "Egfr", "MuSCs", 3, 1, 0.8489523333333334, 2.546857, 0.6346554277922418, 0.634655427792242, 3, 1, 7.591029666666667, 22.773089, 0.1245742015128762, 0.1245742015128762, 6.444422347919223, 57.999801131273, 0.07906169315303137, 0.0790616931530314),
    @("FAPs", "Ereg", "Egfr", "Neutrophils", 3, 1, 0.8489523333333334, 2.546857, 0.6346554277922418, 0.634655427792242, 3, 1, 4.780457, 14.341371, 0.07845070297336118, 0.07845070297336117, 4.058380124549667, 36.525421120947, 0.04978916445616064, 0.04978916445616065),
    @("FAPs", "Ereg", "Egfr", "Resolving-Mac", 3, 1, 0.8489523333333334, 2.546857, 0.6346554277922418, 0.634655427792242, 2, 0.6666666666666666, 0.172647, 0.517941, 0.002833260191701732, 0.002833260191701732, 0.146569073493, 1.319121661437, 0.001798143959011192, 0.001798143959011192),
    @("MuSCs", "Ereg", "Egfr", "ECs", 3, 1, 0.4686963333333334, 1.406089, 0.3503855991164661, 0.3503855991164662, 3, 1, 1.635346666666667, 4.90604, 0.02683720313876748, 0.02683720313876747, 0.7664809863955556, 6.898328877560001, 0.009403369500387348, 0.00940336950038735),
    @("MuSCs", "Ereg", "Egfr", "FAPs", 3, 1, 0.4686963333333334, 1.406089, 0.3503855991164661, 0.3503855991164662, 3, 1, 45.91636366666668, 137.749091, 0.7535202194331003, 0.7535202194331002, 21.52083129056656, 193.687481615099, 0.2640226335324379, 0.2640226335324379),
    @("MuSCs", "Ereg", "Egfr", "Inflammatory-Mac", 3, 1, 0.4686963333333334, 1.406089, 0.3503855991164661, 0.3503855991164662, 3, 1, 0.8399643333333332, 2.519893, 0.01378441275019327, 0.01378441275019327, 0.3936882031641111, 3.543193828477, 0.004829859719945122, 0.004829859719945123),
    @("MuSCs", "Ereg", "Egfr", "MuSCs", 3, 1, 0.4686963333333334, 1.406089, 0.3503855991164661, 0.3503855991164662, 3, 1, 7.591029666666667, 22.773089, 0.1245742015128762, 0.1245742015128762, 3.557887770991222, 32.02098993892101, 0.0436490062315445, 0.04364900623154451),
    @("MuSCs", "Ereg", "Egfr", "Neutrophils", 3, 1, 0.4686963333333334, 1.406089, 0.3503855991164661, 0.3503855991164662, 3, 1, 4.780457, 14.341371, 0.07845070297336118, 0.07845070297336117, 2.240582667557667, 20.165244008019, 0.02748799656242909, 0.02748799656242909),
    @("MuSCs", "Ereg", "Egfr", "Resolving-Mac", 3, 1, 0.4686963333333334, 1.406089, 0.3503855991164661, 0.3503855991164662, 2, 0.6666666666666666, 0.172647, 0.517941, 0.002833260191701732, 0.002833260191701732, 0.080919015861, 0.7282711427490001, 0.000992733569722245, 0.000992733569722245),
    @("Resolving-Mac", "Ereg", "Egfr", "ECs", 1, 0.3333333333333333, 0.02001, 0.06003, 0.01495897309129185, 0.01495897309129185, 3, 1, 1.635346666666667, 4.90604, 0.02683720313876748, 0.02683720313876747, 0.0327232868, 0.2945095812, 0.0004014569995983558, 0.0004014569995983558),
    @("Resolving-Mac", "Ereg", "Egfr", "FAPs", 1, 0.3333333333333333, 0.02001, 0.06003, 0.01495897309129185, 0.01495897309129185, 3, 1, 45.91636366666668, 137.749091, 0.7535202194331003, 0.7535202194331002, 0.9187864369700002, 8.269077932730001, 0.01127188868624407, 0.01127188868624408),
    @("Resolving-Mac", "Ereg", "Egfr", "Inflammatory-Mac", 1, 0.3333333333333333, 0.02001, 0.06003, 0.01495897309129185, 0.01495897309129185, 3, 1, 0.8399643333333332, 2.519893, 0.01378441275019327, 0.01378441275019327, 0.01680768631, 0.15126917679, 0.0002062006594094013, 0.0002062006594094013),
    @("Resolving-Mac", "Ereg", "Egfr", "MuSCs", 1, 0.3333333333333333, 0.02001, 0.06003, 0.01495897309129185, 0.01495897309129185, 3, 1, 7.591029666666667, 22.773089, 0.1245742015128762, 0.1245742015128762, 0.15189650363, 1.36706853267, 0.001863502128300283, 0.001863502128300283),
    @("Resolving-Mac", "Ereg", "Egfr", "Neutrophils", 1, 0.3333333333333333, 0.02001, 0.06003, 0.01495897309129185, 0.01495897309129185, 3, 1, 4.780457, 14.341371, 0.07845070297336118, 0.07845070297336117, 0.09565694457, 0.8609125011300001, 0.001173541954771439, 0.001173541954771439),
    @("Resolving-Mac", "Ereg", "Egfr", "Resolving-Mac", 1, 0.3333333333333333, 0.02001, 0.06003, 0.01495897309129185, 0.01495897309129185, 2, 0.6666666666666666, 0.172647, 0.517941, 0.002833260191701732, 0.002833260191701732, 0.00345466647, 0.03109199823, 0.00004238266296829458, 0.00004238266296829459)
)

for ($r = 0; $r -lt $table.Length; $r++) {
    $rowData = $table[$r]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowData[$c]
    }
}

